# Update the "income/quarterly" statement workbook: new publication-date
# headers (row 9) plus the refreshed year-to-date figures in column M
# that shifted when the new quarter's report was added.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 9: publication-date headers for columns I and M
$ws.Range("I9").Value = "1402-03-13 (10)"
$ws.Range("M9").Value = "1402-03-13 (2)"

# Column M (12-month cumulative) figure updates
$ws.Range("M14").Value = -26784
$ws.Range("M17").Value = 36510
$ws.Range("M18").Value = -22161
$ws.Range("M20").Value = 91996
$ws.Range("M21").Value = -2104
$ws.Range("M22").Value = 89893
$ws.Range("M24").Value = 89893
